# Bulk-update the "Tested" status column (E) on Sheet1, and move the
# active selection/scroll position, matching the upstream commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column E ("Tested") status updates -----------------------------
# Most rows become "DONE"; a few special rows get distinct markers.
# Cells that previously had no value at all need their style (s="1",
# matching the rest of the row) copied over from column D before the
# value is written, since a bare .Value assignment on a brand-new cell
# leaves it with no style id.

function Set-TestedStatus($firstRow, $lastRow, $status) {
    if ($firstRow -eq $lastRow) {
        $srcRange = "D$firstRow"
        $dstRange = "E$firstRow"
    } else {
        $srcRange = "D$firstRow`:D$lastRow"
        $dstRange = "E$firstRow`:E$lastRow"
    }
    $ws.Range($srcRange).Copy()
    $ws.Range($dstRange).PasteSpecial(-4122)
    $ws.Range($dstRange).Value = $status
}

Set-TestedStatus 2 10 "DONE"
Set-TestedStatus 12 14 "DONE"
Set-TestedStatus 16 17 "DONE"
Set-TestedStatus 19 63 "DONE"
Set-TestedStatus 64 64 " "
Set-TestedStatus 65 77 "DONE"
Set-TestedStatus 79 80 "??"
Set-TestedStatus 81 87 "DONE"
Set-TestedStatus 88 89 "LATER"

# Row 18 loses its "Tested" value entirely (cell stays, but empty).
$ws.Range("E18").ClearContents()

# --- View state -------------------------------------------------------
# Scroll so row 4 is at the top of the viewport, and select the entire
# 11th row (as if the row header had been clicked).
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Rows.Item(11).Select()
